$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reset the "Percent" number-format style on column E (the
# "% of nodes required to take over network" column) back to Normal,
# since the sourced data below is non-numeric ("N/A" placeholders).
for ($r = 2; $r -le 10; $r++) {
    $ws.Range("E$r").Style = "Normal"
}

# Row 2 - Proof of Work
$ws.Range("B2").Value = 7
$ws.Range("C2").Value = "N/A"
$ws.Range("D2").Value = "N/A"
$ws.Range("E2").Value = "N/A"
$ws.Range("F2").Value = "N/A"
$ws.Range("G2").Value = "N/A"

# Row 3 - Proof of Stake
$ws.Range("B3").Value = 30
$ws.Range("C3").Value = "N/A"
$ws.Range("D3").Value = "N/A"
$ws.Range("E3").Value = "N/A"
$ws.Range("F3").Value = "N/A"
$ws.Range("G3").Value = "N/A"

# Row 4 - Delegated Proof of Stake
$ws.Range("B4").Value = "N/A"
$ws.Range("C4").Value = "N/A"
$ws.Range("D4").Value = "N/A"
$ws.Range("E4").Value = "N/A"
$ws.Range("F4").Value = "N/A"
$ws.Range("G4").Value = "N/A"

# Row 5 - Proof of History
$ws.Range("B5").Value = "N/A"
$ws.Range("C5").Value = "N/A"
$ws.Range("D5").Value = "N/A"
$ws.Range("E5").Value = "N/A"
$ws.Range("F5").Value = "N/A"
$ws.Range("G5").Value = "N/A"

# Row 6 - Proof of Stake with Byzantine Fault Tolerance
$ws.Range("B6").Value = "N/A"
$ws.Range("C6").Value = "N/A"
$ws.Range("D6").Value = "N/A"
$ws.Range("E6").Value = "N/A"
$ws.Range("F6").Value = "N/A"
$ws.Range("G6").Value = "N/A"

# Row 7 - Proof of History with Proof of Stake
$ws.Range("B7").Value = "N/A"
$ws.Range("C7").Value = "N/A"
$ws.Range("D7").Value = "N/A"
$ws.Range("E7").Value = "N/A"
$ws.Range("F7").Value = "N/A"
$ws.Range("G7").Value = "N/A"

# Row 8 - zk-proof
$ws.Range("B8").Value = "N/A"
$ws.Range("C8").Value = "N/A"
$ws.Range("D8").Value = "N/A"
$ws.Range("E8").Value = "N/A"
$ws.Range("F8").Value = "N/A"
$ws.Range("G8").Value = "N/A"

# Row 9 - Sharding
$ws.Range("B9").Value = "N/A"
$ws.Range("C9").Value = "N/A"
$ws.Range("D9").Value = "N/A"
$ws.Range("E9").Value = "N/A"
$ws.Range("F9").Value = "N/A"
$ws.Range("G9").Value = "at risk from a 1% attack "

# Row 10 - DAGs
$ws.Range("B10").Value = "N/A"
$ws.Range("C10").Value = "N/A"
$ws.Range("D10").Value = "N/A"
$ws.Range("E10").Value = "N/A"
$ws.Range("F10").Value = "N/A"
$ws.Range("G10").Value = "N/A"

# The "Percent" cell style is no longer used anywhere in the sheet, so remove
# it from the workbook's style list.
$wb.Styles.Item("Percent").Delete()

# Update the saved selection to match where the author last left the cursor.
$ws.Range("G9").Select() | Out-Null
